# Fixes a parsing bug in the reference-dumping tool: the abstract/authors text for
# row 13 ("Impact of lockdown on COVID-19 prevalence and mortality...") still carried
# raw id="ParN"> anchor artifacts and inconsistent author-list spacing left over from
# earlier partial-clean passes. Re-write D13 (Abstract) and E13 (Authors) with the
# fully-cleaned text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cleanedAbstract = @'
Background
This study aimed to assess the impact of 15 days before, 15 days during, and 15 days after the lockdown on the trends in the prevalence and mortality in 27 countries during COVID-19 pandemic.
Methods
Twenty-seven countries were randomly selected from the different continents.
 The information on the trends in the prevalence and mortality due to COVID-19 pandemic in 27 countries was obtained from World Health Organization and lockdown data were obtained from concerned countries and their ministries.
 The impact of lockdown for 15 days before, 15 days during, and 15 days after the lockdown on the prevalence and mortality due to the COVID-19 pandemic in 27 countries was analyzed.
Results
The findings showed that 15 days after the lockdown there was a trend toward a decline, but no significant decline in the mean prevalence and mean mortality rate due to the COVID-19 pandemic compared to 15 days before, and 15 days during the lockdown in 27 countries.
 The mean growth factor for number of cases was 1.18 and for mortality rate was 1.16.
Conclusions
The findings indicate that 15 days after the lockdown, daily cases of COVID-19 and the growth factor of the disease showed a declined trend, but there was no significant decline in the prevalence and mortality.

'@

$cleanedAuthors = @'
[Sultan Ayoub%Meo%sultanmeo@hotmail.com%1,       Abdulelah Adnan%Abukhalaf%Abdulelahabukhalaf@gmail.com%1,       Ali Abdullah%Alomar%AliAlomarMD@gmail.com%1,       Faris Jamal%AlMutairi%faris11300@gmail.com%1,       Adnan Mahmood%Usmani%adnan_mahmood71@hotmail.com%1,       David C.%Klonoff%dklonoff@diabetestechnology.org%1]
'@

$ws.Range("D13").Value2 = $cleanedAbstract
$ws.Range("E13").Value2 = $cleanedAuthors

